# Condensamento dos crimes / Criação do índice mensal de criminalidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Condense the 12 "Roubo *" rows (rows 10-21) into a single "Roubos" row (row 10)
$ws.Range("A10").Value = "Roubos"
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 10

# Condense the 3 "Furto *" rows (rows 22-24) into a single "Furtos" row (row 11)
$ws.Range("A11").Value = "Furtos"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 4

# Remove the now-redundant detail rows (old rows 12-24): the remaining Roubo
# rows plus the three Furto rows. Everything below (Sequestro ... Ameaças)
# shifts up to close the gap.
$ws.Range("A12:F24").EntireRow.Delete()

# Adjust the sheet view to the post-edit selection/scroll position
$ws.Range("A11").Select()
